# Adding the changes we made on may 9th
#
# The sheet originally held 20 data rows (rows 2-21, columns A:C).
# This update prepends 7 new data rows (now rows 2-8), shifts the
# previous 20 rows down to rows 9-28, and appends 3 more new rows
# at the end (rows 29-31), bringing the sheet to 30 data rows total
# (A1:C31, including the x/y/z header in row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0.0001527163112768
$ws.Range("B2").Value = -0.024892758578062
$ws.Range("C2").Value = 0.0048869219608604
$ws.Range("A3").Value = -0.0087048299610614
$ws.Range("B3").Value = -0.0183259565383195
$ws.Range("C3").Value = 0.0035124751739203
$ws.Range("A4").Value = -0.0713185146450996
$ws.Range("B4").Value = -0.1441642045974731
$ws.Range("C4").Value = 0.0519235469400882
$ws.Range("A5").Value = -0.09269879758358
$ws.Range("B5").Value = 0.0145080499351024
$ws.Range("C5").Value = -0.007177666760981
$ws.Range("A6").Value = -0.090408056974411
$ws.Range("B6").Value = -0.0394008085131645
$ws.Range("C6").Value = -0.0087048299610614
$ws.Range("A7").Value = -0.07376197725534429
$ws.Range("B7").Value = -0.102472648024559
$ws.Range("C7").Value = 0.0461203269660472
$ws.Range("A8").Value = -0.1244637966156005
$ws.Range("B8").Value = -0.4952589869499206
$ws.Range("C8").Value = 0.25641068816185
$ws.Range("A9").Value = 0.955545961856842
$ws.Range("B9").Value = 0.418595403432846
$ws.Range("C9").Value = 0.5012149214744568
$ws.Range("A10").Value = 0.2005165219306945
$ws.Range("B10").Value = 3.583182811737061
$ws.Range("C10").Value = 0.4506658315658569
$ws.Range("A11").Value = 0.0250454749912023
$ws.Range("B11").Value = 1.217149019241333
$ws.Range("C11").Value = 0.2982549667358398
$ws.Range("A12").Value = 0.4430300295352936
$ws.Range("B12").Value = 1.049771904945374
$ws.Range("C12").Value = -0.152105450630188
$ws.Range("A13").Value = -0.1605048477649688
$ws.Range("B13").Value = 0.6899722814559937
$ws.Range("C13").Value = -0.2005165219306945
$ws.Range("A14").Value = 0.2364048510789871
$ws.Range("B14").Value = -1.487151384353638
$ws.Range("C14").Value = -0.1020144969224929
$ws.Range("A15").Value = -1.693776607513428
$ws.Range("B15").Value = -5.596747398376465
$ws.Range("C15").Value = 1.112385630607605
$ws.Range("A16").Value = -1.431409955024719
$ws.Range("B16").Value = -1.327410221099854
$ws.Range("C16").Value = 1.870163917541504
$ws.Range("A17").Value = 0.08170322328805921
$ws.Range("B17").Value = -2.866485118865967
$ws.Range("C17").Value = 0.2145664244890213
$ws.Range("A18").Value = 0.2112066596746444
$ws.Range("B18").Value = -2.239890098571777
$ws.Range("C18").Value = -0.6490443348884583
$ws.Range("A19").Value = -0.4173736870288849
$ws.Range("B19").Value = -0.5337435007095337
$ws.Range("C19").Value = -0.0603229440748691
$ws.Range("A20").Value = 0.2202169150114059
$ws.Range("B20").Value = 1.032209515571594
$ws.Range("C20").Value = 0.3608686327934265
$ws.Range("A21").Value = 2.595261096954346
$ws.Range("B21").Value = 7.212944030761719
$ws.Range("C21").Value = -0.4940372705459595
$ws.Range("A22").Value = 0.2443460971117019
$ws.Range("B22").Value = 2.484389066696167
$ws.Range("C22").Value = 0.5590944290161133
$ws.Range("A23").Value = 0.3139847218990326
$ws.Range("B23").Value = 0.8843801617622375
$ws.Range("C23").Value = -0.1221730485558509
$ws.Range("A24").Value = -1.060767531394958
$ws.Range("B24").Value = 1.06214189529419
$ws.Range("C24").Value = 0.0241291765123605
$ws.Range("A25").Value = -0.1357648074626922
$ws.Range("B25").Value = -0.5329799056053162
$ws.Range("C25").Value = -0.2561052441596985
$ws.Range("A26").Value = -0.5283984541893005
$ws.Range("B26").Value = -2.495231866836548
$ws.Range("C26").Value = 0.1453859210014343
$ws.Range("A27").Value = -0.8752171993255615
$ws.Range("B27").Value = -1.659415483474731
$ws.Range("C27").Value = -0.0487165041267871
$ws.Range("A28").Value = -0.0806342139840126
$ws.Range("B28").Value = -0.5641340613365173
$ws.Range("C28").Value = 0.087353728711605
$ws.Range("A29").Value = -0.1470658034086227
$ws.Range("B29").Value = -0.0910189226269722
$ws.Range("C29").Value = -0.2229658216238021
$ws.Range("A30").Value = 0.4641048610210418
$ws.Range("B30").Value = 0.3608686327934265
$ws.Range("C30").Value = 0.3602577745914459
$ws.Range("A31").Value = 0.2412917762994766
$ws.Range("B31").Value = 0.2144137024879455
$ws.Range("C31").Value = -0.0186313893646001
